$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Navigation block: R9:T11 (3x3 grid of nav values) ---
$ws.Range("R9").Value = 0
$ws.Range("S9").Value = 2
$ws.Range("T9").Value = 3

$ws.Range("R10").Value = 4
$ws.Range("S10").Value = 5
$ws.Range("T10").Value = 6

$ws.Range("R11").Value = 7
$ws.Range("S11").Value = 8
$ws.Range("T11").Value = 1

# --- Heatmap implementation block: F20:J24 ---
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 3
$ws.Range("I20").Value = 3
$ws.Range("J20").Value = 3

$ws.Range("F21").Value = 3
$ws.Range("G21").Value = 2
$ws.Range("H21").Value = 2
$ws.Range("I21").Value = 2
$ws.Range("J21").Value = 3

$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 1
$ws.Range("I22").Value = 2
$ws.Range("J22").Value = 3

$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 2
$ws.Range("I23").Value = 2
$ws.Range("J23").Value = 3

$ws.Range("F24").Value = 3
$ws.Range("G24").Value = 3
$ws.Range("H24").Value = 3
$ws.Range("I24").Value = 3
$ws.Range("J24").Value = 3

# --- Row 26: index/legend row ---
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 2
$ws.Range("I26").Value = 3
$ws.Range("J26").Value = 4
$ws.Range("K26").Value = 5
$ws.Range("L26").Value = 6
$ws.Range("M26").Value = 7
$ws.Range("N26").Value = 8

# --- Selection / active cell update to match the final authored state ---
$ws.Range("R9:T11").Select() | Out-Null
